$d = $word.ActiveDocument

$replacements = @(
    @{old = "953×7="; new = "494×9="},
    @{old = "409×3="; new = "335×4="},
    @{old = "144×2="; new = "734×7="},
    @{old = "849×4="; new = "179×5="},
    @{old = "398×3="; new = "190×6="},
    @{old = "181×2="; new = "428×8="},
    @{old = "220×3="; new = "252×5="},
    @{old = "203×4="; new = "514×9="},
    @{old = "401×9="; new = "901×7="},
    @{old = "155×4="; new = "311×2="},
    @{old = "308×7="; new = "772×6="},
    @{old = "482×5="; new = "686×8="},
    @{old = "477×6="; new = "622×6="},
    @{old = "322×3="; new = "444×8="},
    @{old = "224×4="; new = "760×6="},
    @{old = "897×7="; new = "930×7="},
    @{old = "306×2="; new = "317×5="},
    @{old = "810×7="; new = "982×5="},
    @{old = "282×5="; new = "863×8="},
    @{old = "508×2="; new = "806×7="},
    @{old = "836×7="; new = "431×9="},
    @{old = "352×9="; new = "546×8="},
    @{old = "383×4="; new = "480×4="},
    @{old = "236×5="; new = "714×6="},
    @{old = "369×5="; new = "201×7="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
